$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 154.5
$ws.Range("I6").Value = 154.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 463.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -351.5
$ws.Range("N6").ClearContents()
$ws.Range("H17").Value = 554.6957
$ws.Range("J17").Value = 554.6957
$ws.Range("L17").Value = 1664.0871
$ws.Range("N17").Value = -2000.0871
$ws.Range("H74").Value = 6593.2593
$ws.Range("I74").Value = 4565.8887
$ws.Range("K74").Value = 4565.8887
$ws.Range("M74").Value = -3629.8887
$ws.Range("H77").Value = 6593.2593
$ws.Range("I77").Value = 4565.8887
$ws.Range("K77").Value = 22829.4435
$ws.Range("M77").Value = -18149.4435
$ws.Range("H98").Value = 3032.225
$ws.Range("I98").Value = 2527.2122
$ws.Range("J98").Value = 5413
$ws.Range("K98").Value = 2527.2122
$ws.Range("L98").Value = 5413
$ws.Range("M98").Value = -1029.2122
$ws.Range("N98").Value = -8409
$ws.Range("H111").Value = 2596.842
$ws.Range("J111").Value = 2966
$ws.Range("L111").Value = 8898
$ws.Range("N111").Value = -15032
$ws.Range("H112").Value = 3061.625
$ws.Range("J112").Value = 3061.625
$ws.Range("L112").Value = 9184.875
$ws.Range("N112").Value = -11400.875
$ws.Range("H113").Value = 6955.548
$ws.Range("I113").Value = 6206.04
$ws.Range("J113").Value = 8057.7646
$ws.Range("K113").Value = 6206.04
$ws.Range("L113").Value = 8057.7646
$ws.Range("M113").Value = -2952.04
$ws.Range("N113").Value = -14565.7646
$ws.Range("H120").Value = 70000
$ws.Range("J120").Value = 70000
$ws.Range("L120").Value = 70000
$ws.Range("N120").Value = -79676
$ws.Range("H122").Value = 3032.225
$ws.Range("I122").Value = 2527.2122
$ws.Range("J122").Value = 5413
$ws.Range("K122").Value = 7581.6366
$ws.Range("L122").Value = 16239
$ws.Range("M122").Value = -5131.6366
$ws.Range("N122").Value = -21139
$ws.Range("H138").Value = 3445.125
$ws.Range("I138").Value = 3985.6924
$ws.Range("J138").Value = 3075.2632
$ws.Range("K138").Value = 11957.0772
$ws.Range("L138").Value = 9225.7896
$ws.Range("M138").Value = -6817.0772
$ws.Range("N138").Value = -19505.7896
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H18").Value = 2000
$ws.Range("I18").Value = 2000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 2000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1678
$ws.Range("N18").ClearContents()
$ws.Range("H45").Value = 2317.5625
$ws.Range("I45").Value = 1321.909
$ws.Range("K45").Value = 1321.909
$ws.Range("M45").Value = -944.9090000000001
$ws.Range("H61").Value = 11065.167
$ws.Range("I61").Value = 7842.357
$ws.Range("K61").Value = 7842.357
$ws.Range("M61").Value = -7630.357
$ws.Range("H97").Value = 4369.647
$ws.Range("J97").Value = 2183.8
$ws.Range("L97").Value = 2183.8
$ws.Range("N97").Value = -3175.8
$ws.Range("H132").Value = 1706
$ws.Range("I132").Value = 1706
$ws.Range("K132").Value = 5118
$ws.Range("M132").Value = -2588
$ws.Range("H136").Value = 11065.167
$ws.Range("I136").Value = 7842.357
$ws.Range("K136").Value = 23527.071
$ws.Range("M136").Value = -20977.071
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2614.0952
$ws.Range("I99").Value = 2595.2104
$ws.Range("K99").Value = 2595.2104
$ws.Range("M99").Value = -1097.2104
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3774.85
$ws.Range("I122").Value = 3467.8125
$ws.Range("K122").Value = 10403.4375
$ws.Range("M122").Value = -7953.4375
$ws.Range("H132").Value = 2935.9092
$ws.Range("I132").Value = 2301.5
$ws.Range("J132").Value = 4627.6665
$ws.Range("K132").Value = 6904.5
$ws.Range("L132").Value = 13882.9995
$ws.Range("M132").Value = -4374.5
$ws.Range("N132").Value = -18942.9995
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2899.4
$ws.Range("I68").Value = 800
$ws.Range("J68").Value = 3424.25
$ws.Range("K68").Value = 2400
$ws.Range("L68").Value = 10272.75
$ws.Range("M68").Value = -1589
$ws.Range("N68").Value = -11894.75
$ws.Range("H71").Value = 2899.4
$ws.Range("I71").Value = 800
$ws.Range("J71").Value = 3424.25
$ws.Range("K71").Value = 7200
$ws.Range("L71").Value = 30818.25
$ws.Range("M71").Value = -3144
$ws.Range("N71").Value = -38930.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H97").Value = 700.5833
$ws.Range("I97").Value = 539.5
$ws.Range("J97").Value = 1506
$ws.Range("K97").Value = 539.5
$ws.Range("L97").Value = 1506
$ws.Range("M97").Value = -43.5
$ws.Range("N97").Value = -2498
$ws.Range("H126").Value = 4963.143
$ws.Range("I126").Value = 2370
$ws.Range("K126").Value = 7110
$ws.Range("M126").Value = -4640
$ws.Range("H132").Value = 2554.3076
$ws.Range("I132").Value = 2600.5
$ws.Range("K132").Value = 7801.5
$ws.Range("M132").Value = -5271.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3305.3333
$ws.Range("J22").Value = 4500.5
$ws.Range("L22").Value = 4500.5
$ws.Range("N22").Value = -5090.5
$ws.Range("H27").Value = 3305.3333
$ws.Range("J27").Value = 4500.5
$ws.Range("L27").Value = 4500.5
$ws.Range("N27").Value = -4714.5
$ws.Range("H36").Value = 68999.5
$ws.Range("J36").Value = 68999.5
$ws.Range("L36").Value = 68999.5
$ws.Range("N36").Value = -70123.5
$ws.Range("H55").Value = 407.25
$ws.Range("I55").Value = 296.33334
$ws.Range("J55").Value = 740
$ws.Range("K55").Value = 296.33334
$ws.Range("L55").Value = 740
$ws.Range("M55").Value = -123.33334
$ws.Range("N55").Value = -1086
$ws.Range("H136").Value = 2709.4473
$ws.Range("I136").Value = 1155.3572
$ws.Range("J136").Value = 3616
$ws.Range("K136").Value = 3466.0716
$ws.Range("L136").Value = 10848
$ws.Range("M136").Value = -916.0715999999998
$ws.Range("N136").Value = -15948
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2113.2727
$ws.Range("I126").Value = 1921.7778
$ws.Range("J126").Value = 2975
$ws.Range("K126").Value = 5765.3334
$ws.Range("L126").Value = 8925
$ws.Range("M126").Value = -3295.3334
$ws.Range("N126").Value = -13865
